$d = $word.ActiveDocument

# Locate the anchor paragraph ("[x] Home shows API base URL + backend
# reachability") which is immediately followed by a blank paragraph and
# then "Artifacts created:". The new checklist items are inserted right
# after the anchor, before that blank paragraph.
$anchorText = "[x] Home shows API base URL + backend reachability"
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pText = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($pText -eq $anchorText) {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Anchor paragraph not found: $anchorText"
}

$newLines = @(
    "[x] Ride + delivery requests use real backend calls by default (mock off)",
    "[x] Delivery status polling implemented (client-side)",
    "[x] Windows backend runs with CommonJS tsconfig (fixes TS2307)"
)

$idx = $anchorIndex
foreach ($line in $newLines) {
    $d.Paragraphs.Item($idx).Range.InsertParagraphAfter()
    $idx = $idx + 1
    $d.Paragraphs.Item($idx).Range.InsertAfter($line)
}
